# Generate Report for Handback
# Renames the two tracked files in the handback-status report and refreshes
# their corresponding generated-xliff names / timestamps, matching a fresh
# handback run (rows 2 & 3 now collapse onto the same generated xlf).

$wb = $excel.ActiveWorkbook

$oldUuid1 = "58dc30ff-4c75-47e5-95bf-406b1a6b723e"
$newUuid1 = "fdaf2da4-42ac-40e7-bfb6-8b313930fd4f"
$oldUuid2 = "5ec2340e-2472-47a6-8ec5-01d3c83f9145"
$newUuid2 = "ffffbeae2fe0-e746-4aac-85c1-f2ed4616f388"

$newHash1 = "7371ffeec18d473d8a0cb191f85ba38c8143aaff"

$overviewDate      = "2016-08-21 03:05:38"
$zhcnHandoffDate   = "2016-08-21 03:05:34"
$zhcnHandbackDate  = "2016-08-21 03:05:50"
$dedeHandoffDate   = "2016-08-21 03:05:38"
$dedeHandbackDate  = "2016-08-21 03:05:57"

$zhcnXlf = "$newUuid1.$newHash1.zh-cn.xlf"
$dedeXlf = "$newUuid1.$newHash1.de-de.xlf"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newUuid1.md"
$wsOverview.Range("B2").Value = "e2e\$newUuid1.md"
$wsOverview.Range("G2").Value = $overviewDate

$wsOverview.Range("A3").Value = "$newUuid2.md"
$wsOverview.Range("B3").Value = "e2e\$newUuid2.md"
$wsOverview.Range("G3").Value = $overviewDate

if ($wsOverview.Hyperlinks.Count() -ge 1) {
    $wsOverview.Hyperlinks.Item(1).TextToDisplay = "e2e\$newUuid1.md"
}
if ($wsOverview.Hyperlinks.Count() -ge 2) {
    $wsOverview.Hyperlinks.Item(2).TextToDisplay = "e2e\$newUuid2.md"
}

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = "$newUuid1.md"
$wsZhCn.Range("G2").Value = $zhcnXlf
$wsZhCn.Range("H2").Value = $zhcnHandoffDate
$wsZhCn.Range("I2").Value = "$newUuid1.md"
$wsZhCn.Range("J2").Value = $zhcnXlf
$wsZhCn.Range("K2").Value = $zhcnHandbackDate

$wsZhCn.Range("A3").Value = "$newUuid2.md"
$wsZhCn.Range("G3").Value = $zhcnXlf
$wsZhCn.Range("H3").Value = $zhcnHandoffDate
$wsZhCn.Range("I3").Value = "$newUuid2.md"
$wsZhCn.Range("J3").Value = $zhcnXlf
$wsZhCn.Range("K3").Value = $zhcnHandbackDate

if ($wsZhCn.Hyperlinks.Count() -ge 1) {
    $wsZhCn.Hyperlinks.Item(1).TextToDisplay = "$newUuid1.md"
}
if ($wsZhCn.Hyperlinks.Count() -ge 2) {
    $wsZhCn.Hyperlinks.Item(2).TextToDisplay = "$newUuid1.md"
}
if ($wsZhCn.Hyperlinks.Count() -ge 3) {
    $wsZhCn.Hyperlinks.Item(3).TextToDisplay = "$newUuid2.md"
}
if ($wsZhCn.Hyperlinks.Count() -ge 4) {
    $wsZhCn.Hyperlinks.Item(4).TextToDisplay = "$newUuid2.md"
}

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = "$newUuid1.md"
$wsDeDe.Range("G2").Value = $dedeXlf
$wsDeDe.Range("H2").Value = $dedeHandoffDate
$wsDeDe.Range("I2").Value = "$newUuid1.md"
$wsDeDe.Range("J2").Value = $dedeXlf
$wsDeDe.Range("K2").Value = $dedeHandbackDate

$wsDeDe.Range("A3").Value = "$newUuid2.md"
$wsDeDe.Range("G3").Value = $dedeXlf
$wsDeDe.Range("H3").Value = $dedeHandoffDate
$wsDeDe.Range("I3").Value = "$newUuid2.md"
$wsDeDe.Range("J3").Value = $dedeXlf
$wsDeDe.Range("K3").Value = $dedeHandbackDate

if ($wsDeDe.Hyperlinks.Count() -ge 1) {
    $wsDeDe.Hyperlinks.Item(1).TextToDisplay = "$newUuid1.md"
}
if ($wsDeDe.Hyperlinks.Count() -ge 2) {
    $wsDeDe.Hyperlinks.Item(2).TextToDisplay = "$newUuid1.md"
}
if ($wsDeDe.Hyperlinks.Count() -ge 3) {
    $wsDeDe.Hyperlinks.Item(3).TextToDisplay = "$newUuid2.md"
}
if ($wsDeDe.Hyperlinks.Count() -ge 4) {
    $wsDeDe.Hyperlinks.Item(4).TextToDisplay = "$newUuid2.md"
}
